# Update AHP pairwise comparison matrix for "Matriz_Económico" and propagate
# the resulting recomputed weights/results into the dependent sheets.
# (The source workbook stores computed values as static numbers - there are
# no formulas - so the derived sheets are updated with their final values.)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matriz_Económico - raw pairwise comparison inputs changed.
# ---------------------------------------------------------------------
$wsMatriz = $wb.Worksheets.Item("Matriz_Económico")

$wsMatriz.Range("D2").Value = 0.3333333333333333
$wsMatriz.Range("F2").Value = 0.3333333333333333
$wsMatriz.Range("G2").Value = 0.3333333333333333
$wsMatriz.Range("M2").Value = 0.3333333333333333
$wsMatriz.Range("N2").Value = 7

$wsMatriz.Range("D3").Value = 0.3333333333333333
$wsMatriz.Range("F3").Value = 0.3333333333333333
$wsMatriz.Range("G3").Value = 0.3333333333333333
$wsMatriz.Range("M3").Value = 0.3333333333333333
$wsMatriz.Range("N3").Value = 7

$wsMatriz.Range("B4").Value = 3
$wsMatriz.Range("C4").Value = 3
$wsMatriz.Range("E4").Value = 3
$wsMatriz.Range("J4").Value = 3
$wsMatriz.Range("K4").Value = 3

$wsMatriz.Range("D5").Value = 0.3333333333333333
$wsMatriz.Range("F5").Value = 0.3333333333333333
$wsMatriz.Range("G5").Value = 0.3333333333333333
$wsMatriz.Range("M5").Value = 0.3333333333333333
$wsMatriz.Range("N5").Value = 7

$wsMatriz.Range("B6").Value = 3
$wsMatriz.Range("C6").Value = 3
$wsMatriz.Range("E6").Value = 3
$wsMatriz.Range("J6").Value = 3
$wsMatriz.Range("K6").Value = 3

$wsMatriz.Range("B7").Value = 3
$wsMatriz.Range("C7").Value = 3
$wsMatriz.Range("E7").Value = 3
$wsMatriz.Range("J7").Value = 3
$wsMatriz.Range("K7").Value = 3

$wsMatriz.Range("D10").Value = 0.3333333333333333
$wsMatriz.Range("F10").Value = 0.3333333333333333
$wsMatriz.Range("G10").Value = 0.3333333333333333
$wsMatriz.Range("M10").Value = 0.3333333333333333
$wsMatriz.Range("N10").Value = 7

$wsMatriz.Range("D11").Value = 0.3333333333333333
$wsMatriz.Range("F11").Value = 0.3333333333333333
$wsMatriz.Range("G11").Value = 0.3333333333333333
$wsMatriz.Range("M11").Value = 0.3333333333333333

$wsMatriz.Range("B13").Value = 3
$wsMatriz.Range("C13").Value = 3
$wsMatriz.Range("E13").Value = 3
$wsMatriz.Range("J13").Value = 3
$wsMatriz.Range("K13").Value = 3

$wsMatriz.Range("B14").Value = 0.1428571428571428
$wsMatriz.Range("C14").Value = 0.1428571428571428
$wsMatriz.Range("E14").Value = 0.1428571428571428
$wsMatriz.Range("J14").Value = 0.1428571428571428

# ---------------------------------------------------------------------
# 2) Pesos_Locales_Económico - recomputed local weights for "Económico".
# ---------------------------------------------------------------------
$wsPesos = $wb.Worksheets.Item("Pesos_Locales_Económico")

$wsPesos.Range("B2").Value = 0.06859393436079969
$wsPesos.Range("B3").Value = 0.06859393436079969
$wsPesos.Range("B4").Value = 0.1451001396860344
$wsPesos.Range("B5").Value = 0.06859393436079965
$wsPesos.Range("B6").Value = 0.1451001396860344
$wsPesos.Range("B7").Value = 0.1451001396860344
$wsPesos.Range("B8").Value = 0.01428744057464852
$wsPesos.Range("B9").Value = 0.01428744057464852
$wsPesos.Range("B10").Value = 0.06859393436079965
$wsPesos.Range("B11").Value = 0.0697866565851307
$wsPesos.Range("B12").Value = 0.008867907233478526
$wsPesos.Range("B13").Value = 0.1451001396860344
$wsPesos.Range("B14").Value = 0.009419377695460395
$wsPesos.Range("B15").Value = 0.01428744057464852
$wsPesos.Range("B16").Value = 0.01428744057464852

# ---------------------------------------------------------------------
# 3) Resultados - recomputed global weights (rows stay alphabetical).
# ---------------------------------------------------------------------
$wsResultados = $wb.Worksheets.Item("Resultados")

$wsResultados.Range("B2").Value = 0.06510823668815049
$wsResultados.Range("B3").Value = 0.03913143912842076
$wsResultados.Range("B4").Value = 0.04188275646433134
$wsResultados.Range("B5").Value = 0.07867139937636231
$wsResultados.Range("B6").Value = 0.0288494383774855
$wsResultados.Range("B7").Value = 0.03011693642419917
$wsResultados.Range("B8").Value = 0.08795772609039255
$wsResultados.Range("B9").Value = 0.0497829620274733
$wsResultados.Range("B10").Value = 0.09778289071127783
$wsResultados.Range("B11").Value = 0.1313740596055269
$wsResultados.Range("B12").Value = 0.1070658018184729
$wsResultados.Range("B13").Value = 0.02832832624258099
$wsResultados.Range("B14").Value = 0.07124288752260619
$wsResultados.Range("B15").Value = 0.0782350071413063
$wsResultados.Range("B16").Value = 0.06447013238141347

# ---------------------------------------------------------------------
# 4) Ranking_Alternativas - re-sorted ranking (descending by weight).
#    Ranks stay 1..15, but two name/weight pairs swap position because
#    the updated weights changed their relative order.
# ---------------------------------------------------------------------
$wsRanking = $wb.Worksheets.Item("Ranking_Alternativas")

$wsRanking.Range("B2").Value = 0.1313740596055269
$wsRanking.Range("B3").Value = 0.1070658018184729
$wsRanking.Range("B4").Value = 0.09778289071127783
$wsRanking.Range("B5").Value = 0.08795772609039255

$wsRanking.Range("A6").Value = "Jean y Marie Thierry"
$wsRanking.Range("B6").Value = 0.07867139937636231

$wsRanking.Range("A7").Value = "Reina Isabel 2"
$wsRanking.Range("B7").Value = 0.0782350071413063

$wsRanking.Range("B8").Value = 0.07124288752260619

$wsRanking.Range("A9").Value = "Baron"
$wsRanking.Range("B9").Value = 0.06510823668815049

$wsRanking.Range("A10").Value = "Rodelillo"
$wsRanking.Range("B10").Value = 0.06447013238141347

$wsRanking.Range("B11").Value = 0.0497829620274733
$wsRanking.Range("B12").Value = 0.04188275646433134
$wsRanking.Range("B13").Value = 0.03913143912842076
$wsRanking.Range("B14").Value = 0.03011693642419917
$wsRanking.Range("B15").Value = 0.0288494383774855
$wsRanking.Range("B16").Value = 0.02832832624258099
